# Auto-applies the cryptos.xlsx price/volume/coin updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '24.853.15'
$ws.Range("E2").Value = '  +0.60%  '

# Row 3
$ws.Range("D3").Value = '1.713.83'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9985'
$ws.Range("E4").Value = '  -0.36%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.32'
$ws.Range("E5").Value = '  +0.72%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9982'
$ws.Range("E6").Value = '  -0.35%  '

# Row 7
$ws.Range("E7").Value = '  -0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4072'
$ws.Range("E8").Value = '  +0.36%  '

# Row 9
$ws.Range("E9").Value = '  -0.59%  '

# Row 10
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.71'
$ws.Range("E10").Value = '  +0.87%  '

# Row 11
$ws.Range("B11").Value = 'BinanceUSD'
$ws.Range("C11").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9986'
$ws.Range("E11").Value = '  -0.42%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08847'
$ws.Range("E12").Value = '  -0.37%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.46'
$ws.Range("E13").Value = '  +11.68%  '

# Row 14
$ws.Range("E14").Value = '  -1.87%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.164'
$ws.Range("E15").Value = '  +0.13%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001365'
$ws.Range("E16").Value = '  +2.96%  '

# Row 17
$ws.Range("D17").Value = '1.716.93'
$ws.Range("E17").Value = '  +1.53%  '

# Row 18
$ws.Range("E18").Value = '  -2.10%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07222'
$ws.Range("E19").Value = '  +1.63%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.79'
$ws.Range("E20").Value = '  +4.74%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.335'
$ws.Range("E21").Value = '  +2.89%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9987'
$ws.Range("E22").Value = '  -0.27%  '

# Row 23
$ws.Range("E23").Value = '  -2.00%  '

# Row 24
$ws.Range("D24").Value = '24.832.88'
$ws.Range("E24").Value = '  +0.49%  '

# Row 25
$ws.Range("E25").Value = '  -3.85%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.338'
$ws.Range("E26").Value = '  -0.52%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.51'
$ws.Range("E27").Value = '  +1.51%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.03'
$ws.Range("E28").Value = '  +2.18%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.965'
$ws.Range("E29").Value = '  +16.30%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.577'
$ws.Range("E30").Value = '  -7.12%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '146.08'
$ws.Range("E31").Value = '  +4.43%  '

# Row 32
$ws.Range("D32").Value = '1.905.22'
$ws.Range("E32").Value = '  +1.53%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.233'
$ws.Range("E33").Value = '  +14.01%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08880'
$ws.Range("E34").Value = '  -1.55%  '

# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.062'
$ws.Range("E35").Value = '  -0.68%  '

# Row 36
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.03167'
$ws.Range("E36").Value = '  +5.83%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.279'
$ws.Range("E37").Value = '  -8.23%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2867'
$ws.Range("E38").Value = '  +3.04%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8533'
$ws.Range("E39").Value = '  +9.67%  '

# Row 40
$ws.Range("E40").Value = '  -1.61%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09270'
$ws.Range("E41").Value = '  +0.11%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.20'
$ws.Range("E42").Value = '  -1.52%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.486'
$ws.Range("E43").Value = '  +1.42%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.66'
$ws.Range("E44").Value = '  +9.25%  '

# Row 45
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7515'
$ws.Range("E45").Value = '  +3.89%  '

# Row 46
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.711'
$ws.Range("E46").Value = '  +2.66%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.276'
$ws.Range("E47").Value = '  +1.56%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.412'
$ws.Range("E48").Value = '  +3.78%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9993'
$ws.Range("E49").Value = '  -0.16%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '140.85'
$ws.Range("E50").Value = '  +0.55%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08312'
$ws.Range("E51").Value = '  +4.02%  '
